# Added Multibill Verification test case
# Updates the result/date cells on the three "CMCAutopay*" sheets from
# Fail -> Pass with fresh execution timestamps.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("CMCAutopayPC_27")
$ws1.Range("A2").Value = "Pass"
$ws1.Range("B2").Value = "Wed Feb 12 15:18:25 IST 2025"

$ws2 = $wb.Worksheets.Item("CMCAutopayCorp_27")
$ws2.Range("A2").Value = "Pass"
$ws2.Range("B2").Value = "Wed Feb 12 15:21:32 IST 2025"

$ws3 = $wb.Worksheets.Item("CMCAutopayPS_27")
$ws3.Range("A2").Value = "Pass"
$ws3.Range("B2").Value = "Wed Feb 12 15:15:11 IST 2025"
